$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Header row (row 1)
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Data rows
# row 2
$ws.Range("B2").Value = "臺灣銀行"
$ws.Range("C2").Value = "活期儲蓄存款"
$ws.Range("D2").Value = "新臺幣"
$ws.Range("E2").Value = "林道春"
$ws.Range("F2").Value = 15049
$ws.Range("G2").Value = "deposit"
$ws.Range("H2").Value = "normal"
$ws.Range("I2").Value = "2011-12-19"
$ws.Range("J2").Value = "柯建銘"
$ws.Range("K2").Value = 629
$ws.Range("L2").Value = "tmp6ad91"
$ws.Range("M2").Value = 41

# row 3
$ws.Range("B3").Value = "臺灣銀行"
$ws.Range("C3").Value = "活期儲蓄仔款"
$ws.Range("D3").Value = "美金"
$ws.Range("E3").Value = "柯建銘"
$ws.Range("F3").Value = 299577.4
$ws.Range("G3").Value = "deposit"
$ws.Range("H3").Value = "normal"
$ws.Range("I3").Value = "2011-12-19"
$ws.Range("J3").Value = "柯建銘"
$ws.Range("K3").Value = 629
$ws.Range("L3").Value = "tmp6ad91"
$ws.Range("M3").Value = 42

# row 4
$ws.Range("B4").Value = "臺灣銀行"
$ws.Range("C4").Value = "支票存款"
$ws.Range("D4").Value = "新臺幣"
$ws.Range("E4").Value = "柯建銘"
$ws.Range("F4").Value = 564
$ws.Range("G4").Value = "deposit"
$ws.Range("H4").Value = "normal"
$ws.Range("I4").Value = "2011-12-19"
$ws.Range("J4").Value = "柯建銘"
$ws.Range("K4").Value = 629
$ws.Range("L4").Value = "tmp6ad91"
$ws.Range("M4").Value = 43

# row 5
$ws.Range("B5").Value = "合作金庫商業銀行"
$ws.Range("C5").Value = "支票存款"
$ws.Range("D5").Value = "新臺幣"
$ws.Range("E5").Value = "林道春"
$ws.Range("F5").Value = 3102
$ws.Range("G5").Value = "deposit"
$ws.Range("H5").Value = "normal"
$ws.Range("I5").Value = "2011-12-19"
$ws.Range("J5").Value = "柯建銘"
$ws.Range("K5").Value = 629
$ws.Range("L5").Value = "tmp6ad91"
$ws.Range("M5").Value = 44

# row 6
$ws.Range("B6").Value = "潼打銀行"
$ws.Range("C6").Value = "活期儲蓄存款"
$ws.Range("D6").Value = "新臺幣"
$ws.Range("E6").Value = "柯建銘"
$ws.Range("F6").Value = 61
$ws.Range("G6").Value = "deposit"
$ws.Range("H6").Value = "normal"
$ws.Range("I6").Value = "2011-12-19"
$ws.Range("J6").Value = "柯建銘"
$ws.Range("K6").Value = 629
$ws.Range("L6").Value = "tmp6ad91"
$ws.Range("M6").Value = 45

# row 7
$ws.Range("B7").Value = "台北富邦商業銀行"
$ws.Range("C7").Value = "支票存款"
$ws.Range("D7").Value = "新臺幣"
$ws.Range("E7").Value = "柯建銘"
$ws.Range("F7").Value = 183
$ws.Range("G7").Value = "deposit"
$ws.Range("H7").Value = "normal"
$ws.Range("I7").Value = "2011-12-19"
$ws.Range("J7").Value = "柯建銘"
$ws.Range("K7").Value = 629
$ws.Range("L7").Value = "tmp6ad91"
$ws.Range("M7").Value = 46

# row 8
$ws.Range("B8").Value = "台北富邦商業銀行"
$ws.Range("C8").Value = "活期存款"
$ws.Range("D8").Value = "新臺幣"
$ws.Range("E8").Value = "柯建銘"
$ws.Range("F8").Value = 111393
$ws.Range("G8").Value = "deposit"
$ws.Range("H8").Value = "normal"
$ws.Range("I8").Value = "2011-12-19"
$ws.Range("J8").Value = "柯建銘"
$ws.Range("K8").Value = 629
$ws.Range("L8").Value = "tmp6ad91"
$ws.Range("M8").Value = 47

# row 9
$ws.Range("B9").Value = "臺灣銀行"
$ws.Range("C9").Value = "活期儲蓄存款"
$ws.Range("D9").Value = "新臺幣"
$ws.Range("E9").Value = "柯建銘"
$ws.Range("F9").Value = 28607
$ws.Range("G9").Value = "deposit"
$ws.Range("H9").Value = "normal"
$ws.Range("I9").Value = "2011-12-19"
$ws.Range("J9").Value = "柯建銘"
$ws.Range("K9").Value = 629
$ws.Range("L9").Value = "tmp6ad91"
$ws.Range("M9").Value = 48

# row 10
$ws.Range("B10").Value = "第一商業銀行"
$ws.Range("C10").Value = "活期儲蓄存款"
$ws.Range("D10").Value = "新臺幣"
$ws.Range("E10").Value = "柯建銘"
$ws.Range("F10").Value = 3330
$ws.Range("G10").Value = "deposit"
$ws.Range("H10").Value = "normal"
$ws.Range("I10").Value = "2011-12-19"
$ws.Range("J10").Value = "柯建銘"
$ws.Range("K10").Value = 629
$ws.Range("L10").Value = "tmp6ad91"
$ws.Range("M10").Value = 49

# row 11
$ws.Range("B11").Value = "第一商業銀行"
$ws.Range("C11").Value = "支票存款"
$ws.Range("D11").Value = "新臺幣"
$ws.Range("E11").Value = "柯建銘"
$ws.Range("F11").Value = 37320
$ws.Range("G11").Value = "deposit"
$ws.Range("H11").Value = "normal"
$ws.Range("I11").Value = "2011-12-19"
$ws.Range("J11").Value = "柯建銘"
$ws.Range("K11").Value = 629
$ws.Range("L11").Value = "tmp6ad91"
$ws.Range("M11").Value = 50

# row 12
$ws.Range("B12").Value = "華南商業銀行"
$ws.Range("C12").Value = "活期儲蓄存款"
$ws.Range("D12").Value = "新臺幣"
$ws.Range("E12").Value = "林道春"
$ws.Range("F12").Value = 824
$ws.Range("G12").Value = "deposit"
$ws.Range("H12").Value = "normal"
$ws.Range("I12").Value = "2011-12-19"
$ws.Range("J12").Value = "柯建銘"
$ws.Range("K12").Value = 629
$ws.Range("L12").Value = "tmp6ad91"
$ws.Range("M12").Value = 51

# row 13
$ws.Range("B13").Value = "玉山商業銀行"
$ws.Range("C13").Value = "活期儲蓄存款"
$ws.Range("D13").Value = "新臺幣"
$ws.Range("E13").Value = "柯建銘"
$ws.Range("F13").Value = 612
$ws.Range("G13").Value = "deposit"
$ws.Range("H13").Value = "normal"
$ws.Range("I13").Value = "2011-12-19"
$ws.Range("J13").Value = "柯建銘"
$ws.Range("K13").Value = 629
$ws.Range("L13").Value = "tmp6ad91"
$ws.Range("M13").Value = 52

# row 14
$ws.Range("B14").Value = "台新國際商業銀行"
$ws.Range("C14").Value = "活期儲蓄存款"
$ws.Range("D14").Value = "新臺幣"
$ws.Range("E14").Value = "柯建銘"
$ws.Range("F14").Value = 1673
$ws.Range("G14").Value = "deposit"
$ws.Range("H14").Value = "normal"
$ws.Range("I14").Value = "2011-12-19"
$ws.Range("J14").Value = "柯建銘"
$ws.Range("K14").Value = 629
$ws.Range("L14").Value = "tmp6ad91"
$ws.Range("M14").Value = 53

# row 15
$ws.Range("B15").Value = "中國信託商業銀行"
$ws.Range("C15").Value = "活期儲蓄存款"
$ws.Range("D15").Value = "新臺幣"
$ws.Range("E15").Value = "柯建銘"
$ws.Range("F15").Value = 1105
$ws.Range("G15").Value = "deposit"
$ws.Range("H15").Value = "normal"
$ws.Range("I15").Value = "2011-12-19"
$ws.Range("J15").Value = "柯建銘"
$ws.Range("K15").Value = 629
$ws.Range("L15").Value = "tmp6ad91"
$ws.Range("M15").Value = 54

# row 16
$ws.Range("B16").Value = "花旗(台灣)商業銀行"
$ws.Range("C16").Value = "活期儲蓄存款"
$ws.Range("D16").Value = "新臺幣"
$ws.Range("E16").Value = "林道春"
$ws.Range("F16").Value = 354
$ws.Range("G16").Value = "deposit"
$ws.Range("H16").Value = "normal"
$ws.Range("I16").Value = "2011-12-19"
$ws.Range("J16").Value = "柯建銘"
$ws.Range("K16").Value = 629
$ws.Range("L16").Value = "tmp6ad91"
$ws.Range("M16").Value = 55

# row 17
$ws.Range("B17").Value = "臺灣中小企業銀行"
$ws.Range("C17").Value = "活期儲蓄存款"
$ws.Range("D17").Value = "新臺幣"
$ws.Range("E17").Value = "柯建銘"
$ws.Range("F17").Value = 16449
$ws.Range("G17").Value = "deposit"
$ws.Range("H17").Value = "normal"
$ws.Range("I17").Value = "2011-12-19"
$ws.Range("J17").Value = "柯建銘"
$ws.Range("K17").Value = 629
$ws.Range("L17").Value = "tmp6ad91"
$ws.Range("M17").Value = 56

# row 18
$ws.Range("B18").Value = "合作金庫商業銀行"
$ws.Range("C18").Value = "活期儲蓄存款"
$ws.Range("D18").Value = "新臺幣"
$ws.Range("E18").Value = "柯建銘"
$ws.Range("F18").Value = 10152
$ws.Range("G18").Value = "deposit"
$ws.Range("H18").Value = "normal"
$ws.Range("I18").Value = "2011-12-19"
$ws.Range("J18").Value = "柯建銘"
$ws.Range("K18").Value = 629
$ws.Range("L18").Value = "tmp6ad91"
$ws.Range("M18").Value = 57

# row 19
$ws.Range("B19").Value = "彰化商業銀行"
$ws.Range("C19").Value = "活期儲蓄存款"
$ws.Range("D19").Value = "新臺幣"
$ws.Range("E19").Value = "柯建銘"
$ws.Range("F19").Value = 118
$ws.Range("G19").Value = "deposit"
$ws.Range("H19").Value = "normal"
$ws.Range("I19").Value = "2011-12-19"
$ws.Range("J19").Value = "柯建銘"
$ws.Range("K19").Value = 629
$ws.Range("L19").Value = "tmp6ad91"
$ws.Range("M19").Value = 58

# row 20
$ws.Range("B20").Value = "中華郵政股份有限公司"
$ws.Range("C20").Value = "活期儲蓄存款."
$ws.Range("D20").Value = "新臺幣"
$ws.Range("E20").Value = "林道春"
$ws.Range("F20").Value = 1891
$ws.Range("G20").Value = "deposit"
$ws.Range("H20").Value = "normal"
$ws.Range("I20").Value = "2011-12-19"
$ws.Range("J20").Value = "柯建銘"
$ws.Range("K20").Value = 629
$ws.Range("L20").Value = "tmp6ad91"
$ws.Range("M20").Value = 59

# row 21
$ws.Range("B21").Value = "國泰世華商業銀行"
$ws.Range("C21").Value = "活期儲蓄存款"
$ws.Range("D21").Value = "新臺幣"
$ws.Range("E21").Value = "柯建銘"
$ws.Range("F21").Value = 533
$ws.Range("G21").Value = "deposit"
$ws.Range("H21").Value = "normal"
$ws.Range("I21").Value = "2011-12-19"
$ws.Range("J21").Value = "柯建銘"
$ws.Range("K21").Value = 629
$ws.Range("L21").Value = "tmp6ad91"
$ws.Range("M21").Value = 60

# row 22
$ws.Range("B22").Value = "渣打國際商業銀行"
$ws.Range("C22").Value = "活期儲蓄存款"
$ws.Range("D22").Value = "新臺幣"
$ws.Range("E22").Value = "柯建銘"
$ws.Range("F22").Value = 31
$ws.Range("G22").Value = "deposit"
$ws.Range("H22").Value = "normal"
$ws.Range("I22").Value = "2011-12-19"
$ws.Range("J22").Value = "柯建銘"
$ws.Range("K22").Value = 629
$ws.Range("L22").Value = "tmp6ad91"
$ws.Range("M22").Value = 61

# row 23
$ws.Range("B23").Value = "合作金庫商業銀行"
$ws.Range("C23").Value = "活期儲蓄存款"
$ws.Range("D23").Value = "新臺幣"
$ws.Range("E23").Value = "林道春"
$ws.Range("F23").Value = 11209
$ws.Range("G23").Value = "deposit"
$ws.Range("H23").Value = "normal"
$ws.Range("I23").Value = "2011-12-19"
$ws.Range("J23").Value = "柯建銘"
$ws.Range("K23").Value = 629
$ws.Range("L23").Value = "tmp6ad91"
$ws.Range("M23").Value = 62
